$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")
$ws.Activate()

# Row 5: Compilation success -> "no", with note "Called wrong method"
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Called wrong method"

# Row 6: Runtime without error -> clear value (was "yes")
$ws.Range("B6").Value = $null

# Row 7: Assertion validity -> clear value and note
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null

# Row 12: Code BLEU score updates
$ws.Range("B12").Value = 0.2879849400220368
$ws.Range("C12").Value = "{'codebleu': 0.2879849400220368, 'ngram_match_score': 0.10031285888349119, 'weighted_ngram_match_score': 0.12854997812773283, 'syntax_match_score': 0.5659340659340659, 'dataflow_match_score': 0.35714285714285715}"

# Update the selected cell on the sheet
$ws.Range("C6").Select()

$wb.Save()
